$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Santa Lucia" / "Timor Oriental" rows (A202/A203)
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# Update the "Datos actualizados" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Agosto de 2020 a las 08:28"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4973741
$ws.Range("C4").Value = 173
$ws.Range("D4").Value = 2540880
$ws.Range("E4").Value = 2271254
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 161607

# Row 36 - Israel
$ws.Range("B36").Value = 78324
$ws.Range("C36").Value = 405
$ws.Range("D36").Value = 53352
$ws.Range("E36").Value = 24407

# Row 56 - Kirguistan
$ws.Range("B56").Value = 38659
$ws.Range("C56").Value = 549
$ws.Range("D56").Value = 30099
$ws.Range("E56").Value = 7113
$ws.Range("G56").Value = 9
$ws.Range("H56").Value = 1447

# Row 57 - Afganistan
$ws.Range("B57").Value = 36896
$ws.Range("C57").Value = 67
$ws.Range("D57").Value = 25840
$ws.Range("E57").Value = 9758
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 1298

# Row 62 - Uzbekistan
$ws.Range("B62").Value = 28069
$ws.Range("C62").Value = 276
$ws.Range("E62").Value = 9113
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 173

# Row 73
$ws.Range("D73").Value = 9157
$ws.Range("E73").Value = 9031
$ws.Range("G73").Value = 15
$ws.Range("H73").Value = 513

# Row 145 - Georgia
$ws.Range("B145").Value = 1206
$ws.Range("C145").Value = 9
$ws.Range("D145").Value = 987
$ws.Range("E145").Value = 202

# Row 166 - Taiwan
$ws.Range("B166").Value = 477
$ws.Range("C166").Value = 1
$ws.Range("E166").Value = 27
